$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (7:50): remove Monday metalografia entry
$ws.Range("B3").Value = "-"

# Row 4 (8:40): shift metalografia marker from Monday to Friday
$ws.Range("B4").Value = "[-, -, -, 'MEC-2B-Metalografia']"
$ws.Range("F4").Value = "[-, -, 'MEC-2B-Metalografia', -]"

# Row 6 (9:50): shift metalografia marker from Monday to Friday
$ws.Range("B6").Value = "[-, -, -, 'MEC-2B-Metalografia']"
$ws.Range("F6").Value = "[-, -, 'MEC-2B-Metalografia', -]"

# Row 7 (10:40): remove Monday metalografia entry
$ws.Range("B7").Value = "-"
